$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new row 12 data: Hours worked (B12) and Completed task description (C12)
$ws.Range("B12").Value = 0.5
$ws.Range("C12").Value = "Setup view and view navigation"

# Update the active selection to C13, matching the saved workbook state
$ws.Range("C13").Select()
